$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. '93.88') need to be
# forced to Text so Excel doesn't silently convert them to floating point
# numbers, matching the source data which stores every value as text.
$textForceCells = @(
    'D5', 'D6', 'D9', 'D10', 'D11', 'D12', 'D16', 'D17', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26', 'D28', 'D29', 'D31', 'D32', 'D35', 'D36', 'D39', 'D40', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48'
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '45.262.05'
$ws.Range('E2').Value = '  -3.61%  '
$ws.Range('D3').Value = '2.438.86'
$ws.Range('E3').Value = '  +7.72%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '293.86'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '93.88'
$ws.Range('E6').Value = '  -6.26%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').Value = '34.36'
$ws.Range('E10').Value = '  -3.64%  '
$ws.Range('D11').Value = '0.0777'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '7.01'
$ws.Range('E12').Value = '  -2.54%  '
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('D14').Value = '2.810.81'
$ws.Range('E14').Value = '  +7.64%  '
$ws.Range('D15').Value = '2.426.97'
$ws.Range('E15').Value = '  +7.03%  '
$ws.Range('D16').Value = '14.18'
$ws.Range('E16').Value = '  +4.30%  '
$ws.Range('D17').Value = '0.839'
$ws.Range('E17').Value = '  +5.57%  '
$ws.Range('D18').Value = '45.283.50'
$ws.Range('E18').Value = '  -3.42%  '
$ws.Range('D19').Value = '12.36'
$ws.Range('E19').Value = '  -4.42%  '
$ws.Range('D20').Value = '0.0₃0938'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  +6.08%  '
$ws.Range('D22').Value = '66.87'
$ws.Range('E22').Value = '  +2.60%  '
$ws.Range('D23').Value = '238.83'
$ws.Range('E23').Value = '  -4.14%  '
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').Value = '1.91'
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '37.12'
$ws.Range('E28').Value = '  -12.29%  '
$ws.Range('D29').Value = '9.60'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('E30').Value = '  +22.44%  '
$ws.Range('D31').Value = '21.45'
$ws.Range('E31').Value = '  +7.89%  '
$ws.Range('D32').Value = '148.97'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0762'
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '2.00'
$ws.Range('E36').Value = '  +17.48%  '
$ws.Range('E37').Value = '  -1.94%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').Value = '14.31'
$ws.Range('E39').Value = '  -12.15%  '
$ws.Range('D40').Value = '3.72'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('D42').Value = '1.998.15'
$ws.Range('E42').Value = '  +12.25%  '
$ws.Range('D43').Value = '3.16'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '88.21'
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('D46').Value = '16.04'
$ws.Range('E46').Value = '  +24.09%  '
$ws.Range('D47').Value = '1.69'
$ws.Range('E47').Value = '  -13.33%  '
$ws.Range('D48').Value = '8.61'
$ws.Range('E48').Value = '  +9.76%  '
$ws.Range('E49').Value = '  +8.60%  '
$ws.Range('D50').Value = '2.677.45'
$ws.Range('E50').Value = '  +7.60%  '
$ws.Range('E51').Value = '  -3.75%  '

# Strip the temporary Text formatting back off so the cells end up with no
# explicit style override, same as in the source workbook.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).ClearFormats()
}
